# Rename the embedded logo pictures' display "Name" metadata.
#
# The document has three inline pictures living in the header/footer
# stories (not the main body):
#   - Section 1 "first page" footer  -> Pearson logo  : image1.png -> image2.png
#   - Section 1 "primary"    footer  -> Pearson logo  : image1.png -> image2.png
#   - Section 1 "first page" header  -> BTec logo      : image2.jpg -> image1.jpg
#
# wdHeaderFooterIndex constants (Word OM):
#   wdHeaderFooterPrimary   = 1
#   wdHeaderFooterFirstPage = 2

$wdHeaderFooterPrimary = 1
$wdHeaderFooterFirstPage = 2

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Select each picture first, then rename it through the live Selection —
# renaming straight off Range.InlineShapes is flaky for footer stories in
# this host, but it is reliable once the shape's range has been selected.

# Pearson Edexcel logo, first-page footer.
$ftFirst = $sec.Footers.Item($wdHeaderFooterFirstPage)
$ftFirst.Range.InlineShapes.Item(1).Range.Select()
$word.Selection.InlineShapes.Item(1).Name = "image2.png"

# Pearson Edexcel logo, primary (default) footer.
$ftPrimary = $sec.Footers.Item($wdHeaderFooterPrimary)
$ftPrimary.Range.InlineShapes.Item(1).Range.Select()
$word.Selection.InlineShapes.Item(1).Name = "image2.png"

# BTec logo, first-page header.
$hdFirst = $sec.Headers.Item($wdHeaderFooterFirstPage)
$hdFirst.Range.InlineShapes.Item(1).Range.Select()
$word.Selection.InlineShapes.Item(1).Name = "image1.jpg"

Write-Output "Renamed logo inline shapes."
